$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2  = @{ D = 44601; J = 270; K = 2200; L = 2500; M = 2350; N = '$/atado 1,5 a 2 kilos'; P = 1175; Q = 2 }
  3  = @{ D = 44789; J = 300; K = 1400; L = 1500; M = 1450; N = '$/atado 1,5 a 2 kilos'; P = 725;  Q = 2 }
  4  = @{ D = 44525; J = 300; K = 1400; L = 1500; M = 1450; N = '$/atado 1,5 a 2 kilos'; P = 725;  Q = 2 }
  5  = @{ D = 44427; J = 250; K = 1300; L = 1500; M = 1400; N = '$/atado 1,5 a 2 kilos'; P = 700;  Q = 2 }
  6  = @{ D = 44253; J = 250; K = 1800; L = 2000; M = 1900; N = '$/atado 1,5 a 2 kilos'; P = 950;  Q = 2 }
  7  = @{ D = 44229; J = 250; K = 1800; L = 2000; M = 1900; N = '$/atado 1,5 a 2 kilos'; P = 950;  Q = 2 }
  8  = @{ D = 44726; J = 250; K = 2500; L = 2800; M = 2650; N = '$/atado 1,5 a 2 kilos'; P = 1325; Q = 2 }
  9  = @{ D = 44616; J = 270; K = 1300; L = 1500; M = 1400; N = '$/atado 1,5 a 2 kilos'; P = 700;  Q = 2 }
  10 = @{ D = 44817; J = 300; K = 900;  L = 1000; M = 950;  N = '$/atado 1,5 a 2 kilos'; P = 475;  Q = 2 }
  11 = @{ D = 44540; J = 300; K = 900;  L = 1000; M = 950;  N = '$/atado 1,5 a 2 kilos'; P = 475;  Q = 2 }
  12 = @{ D = 44365; J = 200; K = 1800; L = 2000; M = 1900; N = '$/atado 1,5 a 2 kilos'; P = 950;  Q = 2 }
  13 = @{ D = 44266; J = 300; K = 1700; L = 1800; M = 1750; N = '$/atado 1,5 a 2 kilos'; P = 875;  Q = 2 }
  14 = @{ D = 44392; J = 250; K = 1800; L = 2000; M = 1900; N = '$/atado 1,5 a 2 kilos'; P = 950;  Q = 2 }
  15 = @{ D = 44435; J = 300; K = 900;  L = 1000; M = 950;  N = '$/atado 1,5 a 2 kilos'; P = 475;  Q = 2 }
  16 = @{ D = 44544; J = 250; K = 900;  L = 1000; M = 950;  N = '$/atado 1,5 a 2 kilos'; P = 475;  Q = 2 }
  17 = @{ D = 44795; J = 250; K = 1800; L = 2000; M = 1900; N = '$/atado 1,5 a 2 kilos'; P = 950;  Q = 2 }
  18 = @{ D = 44302; J = 300; K = 900;  L = 1000; M = 950;  N = '$/atado 1,5 a 2 kilos'; P = 475;  Q = 2 }
  19 = @{ D = 44390; J = 250; K = 2400; L = 2500; M = 2450; N = '$/atado 1,5 a 2 kilos'; P = 1225; Q = 2 }
  20 = @{ D = 44161; J = 270; K = 900;  L = 1000; M = 950;  N = '$/atado 1,5 a 2 kilos'; P = 475;  Q = 2 }
  21 = @{ D = 44243; J = 250; K = 1200; L = 1300; M = 1250; N = '$/atado 1,5 a 2 kilos'; P = 625;  Q = 2 }
  22 = @{ D = 44202; J = 250; K = 1800; L = 2000; M = 1900; N = '$/atado 1,5 a 2 kilos'; P = 950;  Q = 2 }
  23 = @{ D = 44572; J = 300; K = 1400; L = 1500; M = 1450; N = '$/atado 1,5 a 2 kilos'; P = 725;  Q = 2 }
  24 = @{ D = 44385; J = 300; K = 2400; L = 2500; M = 2450; N = '$/atado 1,5 a 2 kilos'; P = 1225; Q = 2 }
  25 = @{ D = 44403; J = 250; K = 1800; L = 2000; M = 1900; N = '$/atado 1,5 a 2 kilos'; P = 950;  Q = 2 }
  26 = @{ D = 44172; J = 200; K = 1300; L = 1500; M = 1400; N = '$/atado 1,5 a 2 kilos'; P = 700;  Q = 2 }
  27 = @{ D = 44257; J = 500; K = 1400; L = 1500; M = 1450; N = '$/atado 1,5 a 2 kilos'; P = 725;  Q = 2 }
  28 = @{ D = 44363; J = 250; K = 2500; L = 2800; M = 2650; N = '$/atado 1,5 a 2 kilos'; P = 1325; Q = 2 }
  29 = @{ D = 44438; J = 300; K = 950;  L = 1000; M = 975;  N = '$/atado 1,5 a 2 kilos'; P = 488;  Q = 2 }
  30 = @{ D = 44181; J = 200; K = 1000; L = 1200; M = 1100; N = '$/atado';               P = 1100; Q = 1 }
  31 = @{ D = 44468; J = 300; K = 900;  L = 1000; M = 950;  N = '$/atado 1,5 a 2 kilos'; P = 475;  Q = 2 }
  32 = @{ D = 44291; J = 250; K = 1800; L = 2000; M = 1900; N = '$/atado 1,5 a 2 kilos'; P = 950;  Q = 2 }
  33 = @{ D = 44447; J = 300; K = 900;  L = 1000; M = 950;  N = '$/atado 1,5 a 2 kilos'; P = 475;  Q = 2 }
}

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  $ws.Range("D$row").Value = $vals.D
  $ws.Range("J$row").Value = $vals.J
  $ws.Range("K$row").Value = $vals.K
  $ws.Range("L$row").Value = $vals.L
  $ws.Range("M$row").Value = $vals.M
  $ws.Range("N$row").Value = $vals.N
  $ws.Range("P$row").Value = $vals.P
  $ws.Range("Q$row").Value = $vals.Q
}
